$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''22.502.37'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = '''1.574.49'
$ws.Range("E3").Value = '  +0.58%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.66%  '

$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").Value = '''290.61'
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").Value = '''0.3701'
$ws.Range("E7").Value = '  -1.10%  '

$ws.Range("D8").Value = '''50.02'
$ws.Range("E8").Value = '  +1.73%  '

$ws.Range("D9").Value = '''0.3395'
$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").Value = '''1.150'
$ws.Range("E10").Value = '  +2.99%  '

$ws.Range("D11").Value = '''0.07569'
$ws.Range("E11").Value = '  +1.82%  '

$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.56%  '

$ws.Range("D13").Value = '''21.26'
$ws.Range("E13").Value = '  +2.75%  '

$ws.Range("D14").Value = '''6.036'
$ws.Range("E14").Value = '  +2.96%  '

$ws.Range("D15").Value = '''7.003'
$ws.Range("E15").Value = '  +2.28%  '

$ws.Range("D16").Value = '''1.572.33'
$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").Value = '''0.00001125'
$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").Value = '''90.53'
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("D19").Value = '''0.06797'
$ws.Range("E19").Value = '  +1.60%  '

$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").Value = '''6.369'
$ws.Range("E21").Value = '  +3.99%  '

$ws.Range("D22").Value = '''16.42'
$ws.Range("E22").Value = '  +1.38%  '

$ws.Range("E23").Value = '  +3.35%  '

$ws.Range("D24").Value = '''22.489.91'
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").Value = '''2.370'
$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").Value = '''2.657'
$ws.Range("E26").Value = '  +5.80%  '

$ws.Range("E27").Value = '  +0.95%  '

$ws.Range("D28").Value = '''149.78'
$ws.Range("E28").Value = '  +2.01%  '

$ws.Range("D29").Value = '''5.065'
$ws.Range("E29").Value = '  +1.52%  '

$ws.Range("D30").Value = '''124.93'
$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("D31").Value = '''1.752.11'
$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("D32").Value = '''1.071'
$ws.Range("E32").Value = '  +9.52%  '

$ws.Range("D33").Value = '''6.225'
$ws.Range("E33").Value = '  +5.69%  '

$ws.Range("D34").Value = '''2.016'
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("D35").Value = '''9.874'
$ws.Range("E35").Value = '  +1.87%  '

$ws.Range("D36").Value = '''0.08411'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").Value = '''0.02486'
$ws.Range("E37").Value = '  +1.54%  '

$ws.Range("D38").Value = '''0.2304'
$ws.Range("E38").Value = '  +2.36%  '

$ws.Range("D39").Value = '''1.349'
$ws.Range("E39").Value = '  -1.96%  '

$ws.Range("D40").Value = '''0.06543'
$ws.Range("E40").Value = '  +3.17%  '

$ws.Range("D41").Value = '''5.446'
$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("D42").Value = '''11.33'
$ws.Range("E42").Value = '  +3.86%  '

$ws.Range("D43").Value = '''0.6254'
$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("D44").Value = '''14.08'
$ws.Range("E44").Value = '  +2.43%  '

$ws.Range("D45").Value = '''1.0000'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").Value = '''3.803'
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("D47").Value = '''0.5887'
$ws.Range("E47").Value = '  +2.58%  '

$ws.Range("D48").Value = '''2.072'
$ws.Range("E48").Value = '  +2.08%  '

$ws.Range("D49").Value = '''127.48'
$ws.Range("E49").Value = '  +2.09%  '

$ws.Range("D50").Value = '''1.239'
$ws.Range("E50").Value = '  +1.22%  '

$ws.Range("D51").Value = '''0.07318'
$ws.Range("E51").Value = '  +0.23%  '
